$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the values in the score columns (C:G and I) for rows 2-13,
# leaving the cell styling/formatting untouched.
$ws.Range("C2:G13").ClearContents()
$ws.Range("I2:I13").ClearContents()

# Update the sheet view: select C2:I13 (active cell C2) and drop the
# previous scrolled/selected state (H2).
$ws.Range("C2:I13").Select()
